$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data row (row 5: "ewrwr fgfsfsdfs", 3234234, 54, 3234234, "Standard Bank")
$ws.Range("A5:E5").ClearContents()

# Add the new data row (row 23), all stored as text values like in the target workbook
$ws.Range("A23").Value = "qeqweqw dasda"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "3123"
$ws.Range("B23").Style = "Normal"

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "132123"
$ws.Range("C23").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2312"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "eqwweq"
